$d = $word.ActiveDocument

# Locate the "Test Execution Summary" table: header row is
# "Total Test Cases | Passed | Failed | Blocked" and the following row
# holds the figures, with the first cell ("Total Test Cases" value)
# currently reading "20". We need that value to become "24", reproducing
# it as two runs ("2" then "4") the same way Word splits a run when you
# replace a single character in the middle of existing text.
$targetTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    if ($tbl.Cell(1,1).Range.Text -like "Total Test Cases*") {
        $targetTable = $tbl
        break
    }
}

$cell = $targetTable.Cell(2,1)
$rng = $cell.Range
# Exclude the trailing cell-mark / paragraph-mark character.
$rng.End = $rng.End - 1

# Remove the "0" in "20" (second character of the cell's text).
$charToRemove = $d.Range($rng.Start + 1, $rng.Start + 2)
$charToRemove.Text = ""

# Insert "4" right after the remaining "2", producing "24" as two runs
# (matching the real Word behaviour of keeping the freshly typed
# character in its own run instead of merging it back into its
# neighbour).
$insertPoint = $d.Range($rng.Start + 1, $rng.Start + 1)
$insertPoint.InsertAfter("4")
$newChar = $d.Range($rng.Start + 1, $rng.Start + 2)
$newChar.Bold = 1
$newChar.Bold = 0
